$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The old sheet had a two-row header (row1 + row2) with merged/odd leftover
# strings ("mation", "pompes)", "Hiver", "Eté", "Année", "(MW)", "(GWh)").
# The new layout collapses this into a single clean header row and shifts
# all the plant data rows up by one.
# ---------------------------------------------------------------------------

# Remove the old second header row; this shifts rows 3..10 up to 2..9.
$ws.Rows(2).Delete()

# --- Build the new single header row (row 1) --------------------------------
# A1:E1 carry plain (unstyled) header text. E1 previously held a styled
# header cell ("Hiver"/etc., style index 1) before the row shift above, so
# its old formatting must be cleared explicitly to end up unstyled like the
# brand-new A1:D1 cells.
$ws.Range("E1").ClearFormats()
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# F1:K1 reuse the existing data font (Arial 9, same as style index 1..3) but
# with only the font applied (no explicit number format). Create a scratch
# named style to coax the engine into emitting that exact xf combination,
# apply it to the header cells, then delete the named style again so the
# workbook's style tables stay clean (only cellXfs grows by one entry).
$tmpStyle = $wb.Styles.Add("__tmp_header_style__")
$tmpStyle.Font.Name = "Arial"
$tmpStyle.Font.Size = 9

$headerRange = $ws.Range("F1:K1")
$headerRange.Style = "__tmp_header_style__"

$wb.Styles.Item("__tmp_header_style__").Delete()

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# --- Selection moves to the first data row ----------------------------------
$ws.Range("A2:K2").Select()
